# Weekly update: add a new Fruta/Plátano price observation for
# "Feria Lagunitas de Puerto Montt" on the row that used to be row 239,
# pushing every later record down by one row (285 -> 286 rows of data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the existing row 239; Excel shifts rows
# 239:285 down to 240:286 and extends the used range accordingly.
$ws.Rows("239:239").Insert()

# Populate the newly inserted row with this week's data point.
$ws.Cells.Item(239, 1).Value  = 4
$ws.Cells.Item(239, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(239, 3).Value  = "Los Lagos"
$ws.Cells.Item(239, 4).Value  = 44522
$ws.Cells.Item(239, 5).Value  = 10
$ws.Cells.Item(239, 6).Value  = "Fruta"
$ws.Cells.Item(239, 7).Value  = 100108
$ws.Cells.Item(239, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(239, 9).Value  = 100108006
$ws.Cells.Item(239, 10).Value = "Plátano"
$ws.Cells.Item(239, 11).Value = "Sin especificar"
$ws.Cells.Item(239, 12).Value = "Primera Pintón"
$ws.Cells.Item(239, 13).Value = 400
$ws.Cells.Item(239, 14).Value = 20000
$ws.Cells.Item(239, 15).Value = 21000
$ws.Cells.Item(239, 16).Value = 20500
$ws.Cells.Item(239, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(239, 18).Value = "Ecuador"
$ws.Cells.Item(239, 19).Value = 1025
$ws.Cells.Item(239, 20).Value = 20
